# PubMedIn-3 regression test refresh:
#  1. Bump the "generated at" timestamp in the footer.
#  2. Add the standard b/i/sub/sup/u character styles (PubMed HTML-tag
#     styling helpers) to the style sheet.

$d = $word.ActiveDocument

# --- 1. Refresh the footer timestamp -------------------------------------
# The date line lives in the primary footer of section 1, not in the main
# document story, so Find/Replace has to run against that Range.
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute(
    "2025-06-30 12:13Z / ", $true, $false, $false, $false, $false,
    $true, 1, $false, "2025-07-02 02:48Z / ", 2) | Out-Null

# --- 2. Add the b / i / sub / sup / u character styles --------------------
$wdStyleTypeCharacter = 2
$wdUnderlineSingle = 1

$styleDefs = @(
    @{ Id = "b";   Prop = "Bold" },
    @{ Id = "i";   Prop = "Italic" },
    @{ Id = "sub"; Prop = "Subscript" },
    @{ Id = "sup"; Prop = "Superscript" },
    @{ Id = "u";   Prop = "Underline" }
)

foreach ($def in $styleDefs) {
    $style = $d.Styles.Add($def.Id, $wdStyleTypeCharacter)
    $style.BaseStyle = "DefaultParagraphFont"
    $style.Priority = 1
    $style.QuickStyle = $true

    switch ($def.Prop) {
        "Bold"        { $style.Font.Bold = $true }
        "Italic"      { $style.Font.Italic = $true }
        "Subscript"   { $style.Font.Subscript = $true }
        "Superscript" { $style.Font.Superscript = $true }
        "Underline"   { $style.Font.Underline = $wdUnderlineSingle }
    }
}
